$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (A1:D1)
$ws.Cells.Item(1,1).Value2 = "mx_state"
$ws.Cells.Item(1,2).Value2 = "mx_municipality"
$ws.Cells.Item(1,3).Value2 = "n_matriculas"
$ws.Cells.Item(1,4).Value2 = "pct_matriculas"

# 2. Title-case the Spanish connector words (de, del, la, las, los, el, y)
#    inside the state (col A) and municipality (col B) name cells, rows 2-1311
#    (row 1312 "TOTAL" label is handled separately below).
for ($i = 2; $i -le 1311; $i++) {
    for ($col = 1; $col -le 2; $col++) {
        $cell = $ws.Cells.Item($i, $col)
        $v = $cell.Value2
        if ($v -ne $null) {
            if ($v -is [string]) {
                $nv = $v -replace '\bdel\b', 'Del'
                $nv = $nv -replace '\bde\b', 'De'
                $nv = $nv -replace '\blas\b', 'Las'
                $nv = $nv -replace '\bla\b', 'La'
                $nv = $nv -replace '\blos\b', 'Los'
                $nv = $nv -replace '\bel\b', 'El'
                $nv = $nv -replace '\by\b', 'Y'
                $cell.Value2 = $nv
            }
        }
    }
}

# 3. Fix the grand-total row label
$ws.Cells.Item(1312,1).Value2 = "Total"

# 4. Remove the trailing footer/notes rows (old rows 1313-1318)
$ws.Range("A1313:A1318").EntireRow.Delete()
